$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value for every cell that changed in this update
# (coin name/link re-ranking plus refreshed price / 1h volume figures).
$updates = [ordered]@{
    "D2" = "310.89"
    "E2" = "-0.68%"
    "D3" = "37.67"
    "E3" = "-0.16%"
    "D4" = "5.175"
    "E4" = "2.07%"
    "D5" = "0.07914"
    "E5" = "1.90%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "4.432"
    "E6" = "1.81%"
    "B7" = "FTXToken"
    "C7" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D7" = "1.923"
    "E7" = "2.16%"
    "B8" = "KuCoinToken"
    "C8" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "D8" = "8.288"
    "E8" = "1.16%"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D9" = "3.001"
    "E9" = "-0.39%"
    "B10" = "MXToken"
    "C10" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D10" = "0.9388"
    "E10" = "2.35%"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D11" = "0.1102"
    "E11" = "-11.46%"
    "B12" = "WazirX"
    "C12" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D12" = "0.1939"
    "E12" = "2.12%"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D13" = "0.09094"
    "E13" = "2.28%"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D14" = "0.03306"
    "E14" = "-2.42%"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D15" = "0.09615"
    "E15" = "-0.92%"
    "B16" = "BitForexToken"
    "C16" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D16" = "0.001379"
    "E16" = "0.80%"
    "B17" = "TigerCash"
    "C17" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D17" = "0.005790"
    "E17" = "-1.07%"
    "B18" = "LEO"
    "C18" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D18" = "3.596"
    "E18" = "1.82%"
    "E19" = "0.03%"
    "D20" = "6.433"
    "E20" = "27.93%"
    "D21" = "0.1282"
    "E21" = "-1.11%"
    "D22" = "0.2522"
    "E22" = "-2.68%"
    "D23" = "0.04407"
    "E23" = "0.46%"
    "E24" = "1.52%"
    "D25" = "0.004615"
    "E25" = "8.82%"
    "D26" = "0.0001361"
    "E26" = "0.72%"
    "D27" = "0.0003994"
    "D39" = "0.02250"
    "E39" = "5.32%"
    "D40" = "0.05110"
    "E40" = "2.81%"
    "D41" = "0.007473"
    "E41" = "-3.49%"
    "D42" = "0.008862"
    "E42" = "-10.41%"
    "E43" = "0.88%"
    "E44" = "3.38%"
    "D45" = "0.009328"
    "E45" = "-3.56%"
    "D46" = "0.00006607"
    "E46" = "1.35%"
    "D47" = "0.00000000751"
    "E47" = "0.03%"
    "D48" = "0.002861"
    "E48" = "-6.84%"
    "D49" = "0.001001"
    "E49" = "-40.73%"
    "D50" = "0.00002102"
    "E50" = "0.03%"
    "D51" = "0.0002002"
    "E51" = "0.03%"
}

# Columns D (Price) and E (Volume(1h)) store numeric-/percent-looking values as
# literal text (t="inlineStr") in the original workbook. If we assign a plain
# numeric-looking string straight to .Value, Excel auto-converts it into a real
# number (and reformats the "%" values), which would not match the source file.
# So, for every D/E cell, force Text format first, write the literal string, then
# clear the formatting again so the cell keeps its original (default) style - this
# runtime does not propagate Range operations across multi-area (union) ranges, so
# each cell is handled individually.
foreach ($ref in $updates.Keys) {
    $col = $ref -replace '[0-9]+$', ''
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($ref).NumberFormat = "@"
    }
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in $updates.Keys) {
    $col = $ref -replace '[0-9]+$', ''
    if ($col -eq "D" -or $col -eq "E") {
        $ws.Range($ref).ClearFormats()
    }
}
